# issue #5: stock data output to json file
# Add a "property_category" column (constant value "stock") to the
# 股票 (stock) sheet, between the existing "total" and "date" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H; this shifts the existing date/legislator_name/
# legislator_id columns (and their formatting) one column to the right.
$ws.Range("H1").EntireColumn.Insert()

# Header for the new column.
$ws.Cells.Item(1, 8).Value = "property_category"

# Every stock row gets the constant category value "stock".
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
